$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the base value in A1; dependent formulas (A4, B4, A5, B5, A6, B6)
# will recalculate automatically.
$ws.Range("A1").Value = 10.456

$excel.CalculateFullRebuild()
